# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Cell values are plain text (prices use "." as thousands separators, e.g. "24.548.69"),
# so numeric-looking entries are forced to Text format before the write, then the
# format is reset to Normal afterwards (matching how the source data is stored).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '24.548.69'
$ws.Range("E2").Value = '  +3.16%  '

# Row 3
$ws.Range("D3").Value = '1.693.53'

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.66%  '

# Row 6
$ws.Range("E6").Value = '  +0.11%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3933'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.30%  '

# Row 8
$ws.Range("E8").Value = '  +1.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.524'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.71%  '

# Row 10
$ws.Range("E10").Value = '  +0.13%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.10'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.31%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08727'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.75%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.202'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.60%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.07%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001312'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.29%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.574'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.93%  '

# Row 17
$ws.Range("D17").Value = '1.689.35'
$ws.Range("E17").Value = '  +1.31%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '99.45'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.16%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07063'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.99%  '

# Row 20
$ws.Range("E20").Value = '  +2.30%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.872'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.26%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.24%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.22%  '

# Row 24
$ws.Range("D24").Value = '24.542.15'
$ws.Range("E24").Value = '  +3.12%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.062'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.82%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.331'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.76%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.26'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.28%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.94'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.98%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.217'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.89%  '

# Row 30
$ws.Range("E30").Value = '  +3.35%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.601'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.30%  '

# Row 32
$ws.Range("D32").Value = '1.876.32'
$ws.Range("E32").Value = '  +1.14%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.086'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.15%  '

# Row 34
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.302'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.83%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08546'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.87%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.23'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.36%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.947'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.37%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2709'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.02%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.99%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02741'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.86%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09019'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.57%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.469'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.58%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7661'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.09%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7158'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.64%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.33'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.45%  '

# Row 46
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.519'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.65%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.198'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.52%  '

# Row 49
$ws.Range("E49").Value = '  +8.47%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '140.51'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.97%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07986'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.60%  '
